# additional F_min expressions tested
$wb = $excel.ActiveWorkbook

# --- Sheet 1: input_stoich_coefficients ---
# PLP -> Ald, T3H -> Hydr, Comp -> Hydrz (header "name" is unchanged)
$ws1 = $wb.Worksheets.Item(1)
$ws1.Range("A1").Value = "Ald"
$ws1.Range("B1").Value = "Hydr"
$ws1.Range("C2").Value = "Hydrz"

# --- Sheet 3: input_concentrations ---
# row2 species headers PLP -> Ald, T3H -> Hydr ("tot" row1 unchanged)
$ws3 = $wb.Worksheets.Item(3)
$ws3.Range("A2").Value = "Ald"
$ws3.Range("B2").Value = "Hydr"

# --- Sheet 4: component_name ---
# T3H -> Hydr
$ws4 = $wb.Worksheets.Item(4)
$ws4.Range("A1").Value = "Hydr"

# --- Sheet 5: heats ---
# Row 2 (B2:P2) numeric values updated slightly (additional F_min expressions tested)
$ws5 = $wb.Worksheets.Item(5)
$ws5.Range("B2").Value = 15.024977
$ws5.Range("C2").Value = 15.049954
$ws5.Range("D2").Value = 15.074930999999999
$ws5.Range("E2").Value = 15.099907999999999
$ws5.Range("F2").Value = 15.124885000000001
$ws5.Range("G2").Value = 15.149862000000001
$ws5.Range("H2").Value = 15.174839
$ws5.Range("I2").Value = 15.199816
$ws5.Range("J2").Value = 15.224793
$ws5.Range("K2").Value = 15.24977
$ws5.Range("L2").Value = 15.274747
$ws5.Range("M2").Value = 15.299723999999999
$ws5.Range("N2").Value = 15.324700999999999
$ws5.Range("O2").Value = 15.349678000000001
$ws5.Range("P2").Value = 15.374655000000001

# --- Sheet 6: targets ---
# Comp -> Hydrz
$ws6 = $wb.Worksheets.Item(6)
$ws6.Range("B1").Value = "Hydrz"

# --- Selection / active-cell updates on each sheet ---
[void]$ws1.Range("C2").Select()

[void]$ws3.Range("E9").Select()

[void]$ws4.Range("E6").Select()

[void]$ws5.Range("N11").Select()

[void]$ws6.Range("F2").Select()

# --- Active sheet switches from input_k_constants_log10 (tab 2) to enthalpies (tab 7) ---
# Re-select sheet 7's existing active cell last so it becomes the active tab (activeTab=6 / tabSelected moves here)
$ws7 = $wb.Worksheets.Item(7)
[void]$ws7.Range("A2").Select()
